$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 11.97902751906419
$ws.Cells.Item(2, 3).Value = 7.84689065220905
$ws.Cells.Item(2, 5).Value = 11.6575429873258
$ws.Cells.Item(2, 6).Value = 16.86991607391233
$ws.Cells.Item(2, 7).Value = 3.661620903701022
$ws.Cells.Item(2, 9).Value = 24.73826245096367
$ws.Cells.Item(2, 11).Value = 9.289709712309877
$ws.Cells.Item(2, 12).Value = 9.785653689101096
$ws.Cells.Item(2, 15).Value = 24.69388541469624

# Row 3
$ws.Cells.Item(3, 2).Value = 11.69010833317379
$ws.Cells.Item(3, 3).Value = 7.831060713654015
$ws.Cells.Item(3, 5).Value = 11.66887114664601
$ws.Cells.Item(3, 6).Value = 15.89584955866808
$ws.Cells.Item(3, 7).Value = 3.663272584161525
$ws.Cells.Item(3, 9).Value = 24.86574163845381
$ws.Cells.Item(3, 11).Value = 9.088736309096127
$ws.Cells.Item(3, 12).Value = 9.770893891731893
$ws.Cells.Item(3, 15).Value = 24.8182894938889

# Row 4
$ws.Cells.Item(4, 2).Value = 11.5104732553419
$ws.Cells.Item(4, 3).Value = 7.82151348027254
$ws.Cells.Item(4, 5).Value = 11.67802608937047
$ws.Cells.Item(4, 6).Value = 15.26997757108491
$ws.Cells.Item(4, 7).Value = 3.664340177004716
$ws.Cells.Item(4, 9).Value = 24.94869804238394
$ws.Cells.Item(4, 11).Value = 8.963941550707185
$ws.Cells.Item(4, 12).Value = 9.763467871891672
$ws.Cells.Item(4, 15).Value = 24.89975329100785

# Row 5
$ws.Cells.Item(5, 2).Value = 11.43681080462263
$ws.Cells.Item(5, 3).Value = 7.817667571401005
$ws.Cells.Item(5, 5).Value = 11.68231011336775
$ws.Cells.Item(5, 6).Value = 15.00819731993403
$ws.Cells.Item(5, 7).Value = 3.66478871230901
$ws.Cells.Item(5, 9).Value = 24.98368215807549
$ws.Cells.Item(5, 11).Value = 8.912801737262779
$ws.Cells.Item(5, 12).Value = 9.760855798463936
$ws.Cells.Item(5, 15).Value = 24.93422730841419

# Row 6
$ws.Cells.Item(6, 2).Value = 11.42455454675701
$ws.Cells.Item(6, 3).Value = 7.817031696819317
$ws.Cells.Item(6, 5).Value = 11.68305489496575
$ws.Cells.Item(6, 6).Value = 14.96433081551593
$ws.Cells.Item(6, 7).Value = 3.664864006774432
$ws.Cells.Item(6, 9).Value = 24.98956246404625
$ws.Cells.Item(6, 11).Value = 8.904294782727398
$ws.Cells.Item(6, 12).Value = 9.760447145278507
$ws.Cells.Item(6, 15).Value = 24.94002880288541

# Row 7
$ws.Cells.Item(7, 2).Value = 11.50948153957967
$ws.Cells.Item(7, 3).Value = 7.821461430427373
$ws.Cells.Item(7, 5).Value = 11.67808162466822
$ws.Cells.Item(7, 6).Value = 15.26647399323137
$ws.Cells.Item(7, 7).Value = 3.664346171464011
$ws.Cells.Item(7, 9).Value = 24.94916507655609
$ws.Cells.Item(7, 11).Value = 8.963252925752087
$ws.Cells.Item(7, 12).Value = 9.763430964733036
$ws.Cells.Item(7, 15).Value = 24.90021305013357

# Row 8
$ws.Cells.Item(8, 2).Value = 11.87993284281078
$ws.Cells.Item(8, 3).Value = 7.841397495175804
$ws.Cells.Item(8, 5).Value = 11.66099264793286
$ws.Cells.Item(8, 6).Value = 16.53996406344768
$ws.Cells.Item(8, 7).Value = 3.662179333413409
$ws.Cells.Item(8, 9).Value = 24.78124575340916
$ws.Cells.Item(8, 11).Value = 9.220742782898277
$ws.Cells.Item(8, 12).Value = 9.780226461260625
$ws.Cells.Item(8, 15).Value = 24.73572554080176

# Row 9
$ws.Cells.Item(9, 2).Value = 12.58436646131448
$ws.Cells.Item(9, 3).Value = 7.881802893766348
$ws.Cells.Item(9, 5).Value = 11.64491571100881
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 3.658352430670344
$ws.Cells.Item(9, 9).Value = 24.4890753368638
$ws.Cells.Item(9, 11).Value = 9.711830910229246
$ws.Cells.Item(9, 12).Value = 9.826024249291935
$ws.Cells.Item(9, 15).Value = 24.45348271349408

# Row 10
$ws.Cells.Item(10, 2).Value = 13.08316922400794
$ws.Cells.Item(10, 3).Value = 7.912211810371559
$ws.Cells.Item(10, 5).Value = 11.64370126724387
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 3.655795593524227
$ws.Cells.Item(10, 9).Value = 24.29698497591296
$ws.Cells.Item(10, 11).Value = 10.06069769744589
$ws.Cells.Item(10, 12).Value = 9.867331757390691
$ws.Cells.Item(10, 15).Value = 24.27071880041977

# Row 11
$ws.Cells.Item(11, 2).Value = 13.30501080661604
$ws.Cells.Item(11, 3).Value = 7.926186101053916
$ws.Cells.Item(11, 5).Value = 11.64543954558996
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 3.654687186836812
$ws.Cells.Item(11, 9).Value = 24.214486378848
$ws.Cells.Item(11, 11).Value = 10.21614257054552
$ws.Cells.Item(11, 12).Value = 9.887742796947521
$ws.Cells.Item(11, 15).Value = 24.19292239626028

# Row 12
$ws.Cells.Item(12, 2).Value = 13.38821233560211
$ws.Cells.Item(12, 3).Value = 7.931496624782272
$ws.Cells.Item(12, 5).Value = 11.64642597105013
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 3.654275287239361
$ws.Cells.Item(12, 9).Value = 24.18394799412464
$ws.Cells.Item(12, 11).Value = 10.27448603647609
$ws.Cells.Item(12, 12).Value = 9.895700712985073
$ws.Cells.Item(12, 15).Value = 24.1642321380559

# Row 13
$ws.Cells.Item(13, 2).Value = 13.37033035600124
$ws.Cells.Item(13, 3).Value = 7.930352098541609
$ws.Cells.Item(13, 5).Value = 11.64619895255333
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 3.654363649541509
$ws.Cells.Item(13, 9).Value = 24.19049376230147
$ws.Cells.Item(13, 11).Value = 10.2619446469296
$ws.Cells.Item(13, 12).Value = 9.893976732522807
$ws.Cells.Item(13, 15).Value = 24.17037685395126

# Row 14
$ws.Cells.Item(14, 2).Value = 13.31187239755796
$ws.Cells.Item(14, 3).Value = 7.926622633824312
$ws.Cells.Item(14, 5).Value = 11.64551412997751
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 3.65465314289886
$ws.Cells.Item(14, 9).Value = 24.21195989908897
$ws.Cells.Item(14, 11).Value = 10.22095323084301
$ws.Cells.Item(14, 12).Value = 9.888392941348851
$ws.Cells.Item(14, 15).Value = 24.19054660184444

# Row 15
$ws.Cells.Item(15, 2).Value = 13.27595820117432
$ws.Cells.Item(15, 3).Value = 7.924340626792446
$ws.Cells.Item(15, 5).Value = 11.64513735494109
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 3.654831484523906
$ws.Cells.Item(15, 9).Value = 24.22519994646918
$ws.Cells.Item(15, 11).Value = 10.19577559167326
$ws.Cells.Item(15, 12).Value = 9.885002360466547
$ws.Cells.Item(15, 15).Value = 24.20300140749831

# Row 16
$ws.Cells.Item(16, 2).Value = 13.0685627239028
$ws.Cells.Item(16, 3).Value = 7.91130126697637
$ws.Cells.Item(16, 5).Value = 11.64363366164764
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 3.655869128230443
$ws.Cells.Item(16, 9).Value = 24.30247470935302
$ws.Cells.Item(16, 11).Value = 10.05046896445214
$ws.Cells.Item(16, 12).Value = 9.866030088905823
$ws.Cells.Item(16, 15).Value = 24.27591060358063

# Row 17
$ws.Cells.Item(17, 2).Value = 12.93997828956824
$ws.Cells.Item(17, 3).Value = 7.903337104095471
$ws.Cells.Item(17, 5).Value = 11.64329709549545
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 3.656519675168298
$ws.Cells.Item(17, 9).Value = 24.35113101164899
$ws.Cells.Item(17, 11).Value = 9.96045594177062
$ws.Cells.Item(17, 12).Value = 9.854803196436047
$ws.Cells.Item(17, 15).Value = 24.32200748662825

# Row 18
$ws.Cells.Item(18, 2).Value = 12.86554724006466
$ws.Cells.Item(18, 3).Value = 7.898769704065018
$ws.Cells.Item(18, 5).Value = 11.64331921225848
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 3.656899004472291
$ws.Cells.Item(18, 9).Value = 24.37957657480952
$ws.Cells.Item(18, 11).Value = 9.908379194175344
$ws.Cells.Item(18, 12).Value = 9.848498542893477
$ws.Cells.Item(18, 15).Value = 24.34902410866674

# Row 19
$ws.Cells.Item(19, 2).Value = 12.84026741544472
$ws.Cells.Item(19, 3).Value = 7.897225604514169
$ws.Cells.Item(19, 5).Value = 11.64336378160654
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 3.657028324849285
$ws.Cells.Item(19, 9).Value = 24.38928672590404
$ws.Cells.Item(19, 11).Value = 9.890696393124784
$ws.Cells.Item(19, 12).Value = 9.846390257556825
$ws.Cells.Item(19, 15).Value = 24.35825781577075

# Row 20
$ws.Cells.Item(20, 2).Value = 12.95371581718333
$ws.Cells.Item(20, 3).Value = 7.904183529511864
$ws.Cells.Item(20, 5).Value = 11.64331060759772
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 3.656449890404925
$ws.Cells.Item(20, 9).Value = 24.34590388576358
$ws.Cells.Item(20, 11).Value = 9.970069794506044
$ws.Cells.Item(20, 12).Value = 9.855982538488044
$ws.Cells.Item(20, 15).Value = 24.31704833712164

# Row 21
$ws.Cells.Item(21, 2).Value = 13.32906534547668
$ws.Cells.Item(21, 3).Value = 7.9277175710737
$ws.Cells.Item(21, 5).Value = 11.64570638309343
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 3.654567899476802
$ws.Cells.Item(21, 9).Value = 24.20563572323425
$ws.Cells.Item(21, 11).Value = 10.23300791475343
$ws.Cells.Item(21, 12).Value = 9.890026864587762
$ws.Cells.Item(21, 15).Value = 24.18460136330798

# Row 22
$ws.Cells.Item(22, 2).Value = 13.56965427066572
$ws.Cells.Item(22, 3).Value = 7.943207071229799
$ws.Cells.Item(22, 5).Value = 11.64918432013021
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 3.653383531439967
$ws.Cells.Item(22, 9).Value = 24.11805432215263
$ws.Cells.Item(22, 11).Value = 10.40179997975187
$ws.Cells.Item(22, 12).Value = 9.913607674887489
$ws.Cells.Item(22, 15).Value = 24.10252571292073

# Row 23
$ws.Cells.Item(23, 2).Value = 13.44170334364057
$ws.Cells.Item(23, 3).Value = 7.934930602699379
$ws.Cells.Item(23, 5).Value = 11.64715357339945
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 3.654011488848891
$ws.Cells.Item(23, 9).Value = 24.16442383483383
$ws.Cells.Item(23, 11).Value = 10.31200802863211
$ws.Cells.Item(23, 12).Value = 9.900901850564946
$ws.Cells.Item(23, 15).Value = 24.14592016634591

# Row 24
$ws.Cells.Item(24, 2).Value = 12.94750665159704
$ws.Cells.Item(24, 3).Value = 7.903800825263224
$ws.Cells.Item(24, 5).Value = 11.64330382708202
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 3.656481423518385
$ws.Cells.Item(24, 9).Value = 24.34826559782223
$ws.Cells.Item(24, 11).Value = 9.965724386426622
$ws.Cells.Item(24, 12).Value = 9.85544889090624
$ws.Cells.Item(24, 15).Value = 24.31928876477154

# Row 25
$ws.Cells.Item(25, 2).Value = 12.39672814024662
$ws.Cells.Item(25, 3).Value = 7.870739592734834
$ws.Cells.Item(25, 5).Value = 11.64740096071654
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 3.65934277670877
$ws.Cells.Item(25, 9).Value = 24.56414793077515
$ws.Cells.Item(25, 11).Value = 9.580834666130029
$ws.Cells.Item(25, 12).Value = 9.81227538779808
$ws.Cells.Item(25, 15).Value = 24.52551870429661
